$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing interface rows 22-27 down by one to make room for the
# new "interface Imortality" entry at row 22. Work bottom-up so we don't
# clobber source rows before they're copied, and only touch the columns
# that actually hold data in each source row (A:J for the fully-populated
# rows, F:G for the last, sparsely-populated row).
$ws.Range("F27:G27").Copy($ws.Range("F28:G28"))
$ws.Range("A26:J26").Copy($ws.Range("A27:J27"))
$ws.Range("A25:J25").Copy($ws.Range("A26:J26"))
$ws.Range("A24:J24").Copy($ws.Range("A25:J25"))
$ws.Range("A23:J23").Copy($ws.Range("A24:J24"))
$ws.Range("A22:J22").Copy($ws.Range("A23:J23"))

# Build the new row 22 by reusing row 21's formatting (blank data cells,
# F holds the styled interface label, G stays blank).
$ws.Range("A21:J21").Copy($ws.Range("A22:J22"))
$ws.Range("G22").ClearContents()
$ws.Range("F22").Value = "interface Imortality"

$ws.Range("G9").Select()
